$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.398.65'
$ws.Range("E2").Value = '  -0.56%  '

$ws.Range("D3").Value = '1.840.58'
$ws.Range("E3").Value = '  -0.84%  '

$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").Value = '''261.04'
$ws.Range("E5").Value = '  -3.88%  '

$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.30%  '

$ws.Range("D7").Value = '''0.5189'
$ws.Range("E7").Value = '  -0.74%  '

$ws.Range("D8").Value = '''0.3266'
$ws.Range("E8").Value = '  -3.29%  '

$ws.Range("D9").Value = '''0.06766'
$ws.Range("E9").Value = '  -0.64%  '

$ws.Range("D10").Value = '''18.62'
$ws.Range("E10").Value = '  -5.75%  '

$ws.Range("D11").Value = '''0.7737'
$ws.Range("E11").Value = '  -2.17%  '

$ws.Range("D12").Value = '''0.07758'
$ws.Range("E12").Value = '  +0.70%  '

$ws.Range("D13").Value = '1.837.53'
$ws.Range("E13").Value = '  -2.26%  '

$ws.Range("D14").Value = '''87.40'
$ws.Range("E14").Value = '  -2.14%  '

$ws.Range("D15").Value = '''4.994'
$ws.Range("E15").Value = '  -2.33%  '

$ws.Range("D16").Value = '''1.000'
$ws.Range("E16").Value = '  -0.53%  '

$ws.Range("D17").Value = '''13.88'
$ws.Range("E17").Value = '  -3.55%  '

$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("D19").Value = '''0.000007926'
$ws.Range("E19").Value = '  -0.29%  '

$ws.Range("D20").Value = '26.435.41'
$ws.Range("E20").Value = '  -0.52%  '

$ws.Range("D21").Value = '2.073.60'
$ws.Range("E21").Value = '  -1.44%  '

$ws.Range("D22").Value = '''4.619'
$ws.Range("E22").Value = '  -1.68%  '

$ws.Range("D23").Value = '''9.528'
$ws.Range("E23").Value = '  -3.97%  '

$ws.Range("D24").Value = '''5.989'
$ws.Range("E24").Value = '  -2.10%  '

$ws.Range("D25").Value = '''145.95'
$ws.Range("E25").Value = '  +0.42%  '

$ws.Range("D26").Value = '''2.184'
$ws.Range("E26").Value = '  -6.65%  '

$ws.Range("D27").Value = '''1.650'
$ws.Range("E27").Value = '  -0.59%  '

$ws.Range("D28").Value = '''16.93'
$ws.Range("E28").Value = '  -1.44%  '

$ws.Range("D29").Value = '''111.31'
$ws.Range("E29").Value = '  -0.90%  '

$ws.Range("D30").Value = '''4.190'
$ws.Range("E30").Value = '  -2.02%  '

$ws.Range("D31").Value = '''4.112'
$ws.Range("E31").Value = '  -3.94%  '

$ws.Range("D32").Value = '''0.08700'
$ws.Range("E32").Value = '  -1.94%  '

$ws.Range("D33").Value = '''0.04820'
$ws.Range("E33").Value = '  -1.61%  '

$ws.Range("D34").Value = '''1.127'
$ws.Range("E34").Value = '  -1.76%  '

$ws.Range("D35").Value = '''0.7175'
$ws.Range("E35").Value = '  -0.31%  '

$ws.Range("D36").Value = '''2.849'
$ws.Range("E36").Value = '  -1.32%  '

$ws.Range("D37").Value = '''3.085'
$ws.Range("E37").Value = '  -4.05%  '

$ws.Range("D38").Value = '''0.01775'
$ws.Range("E38").Value = '  -3.31%  '

$ws.Range("D39").Value = '''2.217'
$ws.Range("E39").Value = '  -3.68%  '

$ws.Range("D40").Value = '''0.4815'
$ws.Range("E40").Value = '  -4.83%  '

$ws.Range("D41").Value = '''111.63'
$ws.Range("E41").Value = '  -2.95%  '

$ws.Range("D42").Value = '''0.8979'
$ws.Range("E42").Value = '  -3.31%  '

$ws.Range("D43").Value = '''6.070'
$ws.Range("E43").Value = '  -0.90%  '

$ws.Range("E44").Value = '  -0.34%  '

$ws.Range("D45").Value = '''7.707'
$ws.Range("E45").Value = '  -3.19%  '

$ws.Range("D46").Value = '''0.05952'
$ws.Range("E46").Value = '  -0.02%  '

$ws.Range("D47").Value = '''0.4148'
$ws.Range("E47").Value = '  -5.44%  '

$ws.Range("D48").Value = '''8.997'
$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("E49").Value = '  -2.65%  '

$ws.Range("E50").Value = '  -7.96%  '

$ws.Range("D51").Value = '''0.8831'
$ws.Range("E51").Value = '  +0.80%  '

